$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header block updates (values unchanged semantically, just confirming
#    text stays the same; VALOR MORA total and counts change below).
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 981673
$ws.Range("C13").Value = 7
$ws.Range("F13").Value = 7

# ---------------------------------------------------------------------------
# 2. Make room for the new worker rows. The sheet currently has:
#      row 16/17 -> two data rows (row17 uses the "last row" bottom-border
#      style), rows 22/23 -> footer signature lines.
#    We need 18 new data rows (18..35) so the footer ends up on rows 40/41.
#    Range(...).Insert() shifts everything below down and fixes up the
#    mergeCells automatically.
# ---------------------------------------------------------------------------
$ws.Range("18:35").Insert()

# Row 17 currently still carries the special "last row" style (bottom
# border etc.) - copy that formatting down onto the new last row (35)
# before we overwrite row 17 with the regular interior-row style.
$ws.Range("B17:J17").Copy()
$ws.Range("B35:J35").PasteSpecial(-4122)

# Apply the regular interior-row style (same as row 16) to rows 17..34.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Fill in the worker rows.
# ---------------------------------------------------------------------------
$rows = @(
  @(17, "1043963433", "MARIA JOSE PARODI CASTILLA", "2501", 18980, 1423500),
  @(18, "1143401919", "LEIDY CATHERINE HERNANDEZ SUAREZ", "2507", 52000, 1300000),
  @(19, "1143401919", "LEIDY CATHERINE HERNANDEZ SUAREZ", "2506", 52000, 1300000),
  @(20, "1143401919", "LEIDY CATHERINE HERNANDEZ SUAREZ", "2505", 52000, 1300000),
  @(21, "1143401919", "LEIDY CATHERINE HERNANDEZ SUAREZ", "2504", 52000, 1300000),
  @(22, "1143401919", "LEIDY CATHERINE HERNANDEZ SUAREZ", "2503", 52000, 1300000),
  @(23, "1143401919", "LEIDY CATHERINE HERNANDEZ SUAREZ", "2502", 52000, 1300000),
  @(24, "1103103779", "VALERIA FELIZZOLA", "2507", 56940, 1423500),
  @(25, "1103103779", "VALERIA FELIZZOLA", "2506", 56940, 1423500),
  @(26, "1103103779", "VALERIA FELIZZOLA", "2505", 56940, 1423500),
  @(27, "1103103779", "VALERIA FELIZZOLA", "2504", 56940, 1423500),
  @(28, "1047504763", "CAMILA CASTILLO HERNANDEZ", "2507", 52000, 1300000),
  @(29, "1047504763", "CAMILA CASTILLO HERNANDEZ", "2506", 52000, 1300000),
  @(30, "1047504763", "CAMILA CASTILLO HERNANDEZ", "2505", 52000, 1300000),
  @(31, "1047504763", "CAMILA CASTILLO HERNANDEZ", "2504", 52000, 1300000),
  @(32, "1047504763", "CAMILA CASTILLO HERNANDEZ", "2503", 52000, 1300000),
  @(33, "1047504763", "CAMILA CASTILLO HERNANDEZ", "2502", 52000, 1300000),
  @(34, "1002249233", "LUSIANA ISABEL MULET BARBOZA", "2501", 38133, 1160000),
  @(35, "1048436591", "MARIA DEL MAR OROZCO BADRAN", "2501", 20800, 1300000)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}

# ---------------------------------------------------------------------------
# 4. Column widths - content got wider (longer names / larger numbers),
#    keep the best-fit columns sized to the new content.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 18.54296875
$ws.Columns.Item(3).ColumnWidth = 16.7265625
$ws.Columns.Item(4).ColumnWidth = 36
$ws.Columns.Item(5).ColumnWidth = 13.54296875
$ws.Columns.Item(6).ColumnWidth = 10.1796875
$ws.Columns.Item(7).ColumnWidth = 14.36328125
$ws.Columns.Item(8).ColumnWidth = 19.36328125
$ws.Columns.Item(9).ColumnWidth = 18.08984375
$ws.Columns.Item(10).ColumnWidth = 15
